# Revenue_Cloud_Complete_Upload_Template.xlsx edit script
# - 11_ProductCatalog: remove the two test rows (Test3 / Testx) that were added,
#   shrinking the sheet from A1:G5 back to A1:G3
# - 13_Product2: widen column C (ProductCode) from 22 to 24 chars and append
#   three bundle-component test products (rows 38-40)
# - 26_ProductCategoryProduct: replace the numeric "Name" column with a proper
#   ProductId column, add a new ProductCategoryId column, and append three
#   new bundle-category rows (8-10)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: 11_ProductCatalog  (dimension A1:G5 -> A1:G3)
# ---------------------------------------------------------------------------
$wsCatalog = $wb.Worksheets.Item("11_ProductCatalog")
$wsCatalog.Rows("4:5").Delete()

# ---------------------------------------------------------------------------
# Sheet: 13_Product2  (dimension A1:R37 -> A1:R40, column C width 22 -> 24)
# ---------------------------------------------------------------------------
$wsProduct = $wb.Worksheets.Item("13_Product2")

# Widen the ProductCode column (index 3) from 22 to 24 characters.
$wsProduct.Columns.Item(3).ColumnWidth = 24 - 0.8333333333333334

# Seed rows 38-40 by copying the formatting of the last existing row (37),
# then overwrite the values/types for each new record.
$wsProduct.Range("A37:R37").Copy()
$wsProduct.Range("A38:R40").PasteSpecial(-4122)

# Row 38 - Test
$wsProduct.Cells.Item(38,1).Value = "01tdp000006tQw5AAE"
$wsProduct.Cells.Item(38,2).Value = "Test"
$wsProduct.Cells.Item(38,3).Value = "TEST"
$wsProduct.Cells.Item(38,5).Value = "Test"
$wsProduct.Cells.Item(38,6).Value = $false
$wsProduct.Cells.Item(38,10).Value = $true
$wsProduct.Cells.Item(38,15).Value = $false

# Row 39 - Manual Test Product
$wsProduct.Cells.Item(39,1).Value = "01tdp000006tdZxAAI"
$wsProduct.Cells.Item(39,2).Value = "Manual Test Product 20250728_141722"
$wsProduct.Cells.Item(39,3).Value = "TEST_MANUAL_1753730242"
$wsProduct.Cells.Item(39,5).Value = "Manual test product"
$wsProduct.Cells.Item(39,6).Value = $true
$wsProduct.Cells.Item(39,10).Value = $true
$wsProduct.Cells.Item(39,15).Value = $false

# Row 40 - Test Product
$wsProduct.Cells.Item(40,1).Value = "01tdp000006tgb3AAA"
$wsProduct.Cells.Item(40,2).Value = "Test Product"
$wsProduct.Cells.Item(40,3).Value = "TEST_PRODUCT"
$wsProduct.Cells.Item(40,5).Value = "kujh"
$wsProduct.Cells.Item(40,6).Value = $true
$wsProduct.Cells.Item(40,10).Value = $true
$wsProduct.Cells.Item(40,15).Value = $false

# ---------------------------------------------------------------------------
# Sheet: 26_ProductCategoryProduct  (dimension A1:B7 -> A1:C10)
# ---------------------------------------------------------------------------
$wsCatProd = $wb.Worksheets.Item("26_ProductCategoryProduct")

# Add column C, inheriting column B's formatting, before changing B's content.
$wsCatProd.Range("B1:B7").Copy()
$wsCatProd.Range("C1:C7").PasteSpecial(-4122)

# Update headers: column B becomes ProductId, column C becomes ProductCategoryId.
$wsCatProd.Cells.Item(1,2).Value = "ProductId"
$wsCatProd.Cells.Item(1,3).Value = "ProductCategoryId"

# Replace the old numeric "Name" values in column B with real Product Ids,
# and populate the new ProductCategoryId column.
$wsCatProd.Cells.Item(2,2).Value = "01tdp000006JEGkAAO"
$wsCatProd.Cells.Item(2,3).Value = "0ZGdp0000000Ax3GAE"

$wsCatProd.Cells.Item(3,2).Value = "01tdp000006JEGjAAO"
$wsCatProd.Cells.Item(3,3).Value = "0ZGdp0000000Ax3GAE"

$wsCatProd.Cells.Item(4,2).Value = "01tdp000006JEGlAAO"
$wsCatProd.Cells.Item(4,3).Value = "0ZGdp0000000Ax3GAE"

$wsCatProd.Cells.Item(5,2).Value = "01tdp000006iLGbAAM"
$wsCatProd.Cells.Item(5,3).Value = "0ZGdp0000000AyfGAE"

$wsCatProd.Cells.Item(6,2).Value = "01tdp000006m0jpAAA"
$wsCatProd.Cells.Item(6,3).Value = "0ZGdp0000000AyfGAE"

$wsCatProd.Cells.Item(7,2).Value = "01tdp000006m14nAAA"
$wsCatProd.Cells.Item(7,3).Value = "0ZGdp0000000AyfGAE"

# Seed rows 8-10 by copying the fully-styled row 7 (now 3 columns wide),
# then set the new values for each appended record.
$wsCatProd.Range("A7:C7").Copy()
$wsCatProd.Range("A8:C10").PasteSpecial(-4122)

$wsCatProd.Cells.Item(8,1).Value = "0ZRdp0000000EFtGAM"
$wsCatProd.Cells.Item(8,2).Value = "01tdp000006tcR0AAI"
$wsCatProd.Cells.Item(8,3).Value = "0ZGdp0000000AqbGAE"

$wsCatProd.Cells.Item(9,1).Value = "0ZRdp0000000EHVGA2"
$wsCatProd.Cells.Item(9,2).Value = "01tdp000006teFtAAI"
$wsCatProd.Cells.Item(9,3).Value = "0ZGdp0000000AqbGAE"

$wsCatProd.Cells.Item(10,1).Value = "0ZRdp0000000EJ7GAM"
$wsCatProd.Cells.Item(10,2).Value = "01tdp000006tfiDAAQ"
$wsCatProd.Cells.Item(10,3).Value = "0ZGdp0000000AqbGAE"

# Widen the new/changed columns (B and C) to 20 characters each.
$wsCatProd.Columns.Item(2).ColumnWidth = 20 - 0.8333333333333334
$wsCatProd.Columns.Item(3).ColumnWidth = 20 - 0.8333333333333334

Write-Output "edit complete"
